# Auto-generated: update TPM-derived values for Hspg2-Ptprs sheet
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 155.8700226666667
$ws.Range("H2").Value = 467.610068
$ws.Range("I2").Value = 0.4627663557222626
$ws.Range("J2").Value = 0.4864916976605717
$ws.Range("M2").Value = 2.005664333333333
$ws.Range("N2").Value = 6.016992999999999
$ws.Range("O2").Value = 0.02976781902817159
$ws.Range("P2").Value = 0.03172257287647481
$ws.Range("Q2").Value = 312.6229450983915
$ws.Range("R2").Value = 2813.606505885524
$ws.Range("S2").Value = 0.01377554512946679
$ws.Range("T2").Value = 0.01543276833283744

# Row 3
$ws.Range("G3").Value = 155.8700226666667
$ws.Range("H3").Value = 467.610068
$ws.Range("I3").Value = 0.4627663557222626
$ws.Range("J3").Value = 0.4864916976605717
$ws.Range("O3").Value = 0.4846964599741412
$ws.Range("P3").Value = 0.5165248673390457
$ws.Range("Q3").Value = 5090.303547346843
$ws.Range("R3").Value = 45812.73192612157
$ws.Range("S3").Value = 0.2243012144137149
$ws.Range("T3").Value = 0.2512850595956739

# Row 4
$ws.Range("G4").Value = 155.8700226666667
$ws.Range("H4").Value = 467.610068
$ws.Range("I4").Value = 0.4627663557222626
$ws.Range("J4").Value = 0.4864916976605717
$ws.Range("M4").Value = 10.495667
$ws.Range("N4").Value = 31.487001
$ws.Range("O4").Value = 0.1557753760903259
$ws.Range("P4").Value = 0.1660046278737794
$ws.Range("Q4").Value = 1635.959853191785
$ws.Range("R4").Value = 14723.63867872607
$ws.Range("S4").Value = 0.07208760310458498
$ws.Range("T4").Value = 0.0807598732338264

# Row 5
$ws.Range("G5").Value = 155.8700226666667
$ws.Range("H5").Value = 467.610068
$ws.Range("I5").Value = 0.4627663557222626
$ws.Range("J5").Value = 0.4864916976605717
$ws.Range("M5").Value = 12.4553565
$ws.Range("N5").Value = 24.910713
$ws.Range("O5").Value = 0.1848608423958749
$ws.Range("P5").Value = 0.1313333601264699
$ws.Range("Q5").Value = 1941.416699976414
$ws.Range("R5").Value = 11648.50019985848
$ws.Range("S5").Value = 0.08554737835128655
$ws.Range("T5").Value = 0.06389258932739357

# Row 6
$ws.Range("G6").Value = 155.8700226666667
$ws.Range("H6").Value = 467.610068
$ws.Range("I6").Value = 0.4627663557222626
$ws.Range("J6").Value = 0.4864916976605717
$ws.Range("M6").Value = 9.762884
$ws.Range("N6").Value = 29.288652
$ws.Range("O6").Value = 0.1448995025114864
$ws.Range("P6").Value = 0.1544145717842301
$ws.Range("Q6").Value = 1521.740950372037
$ws.Range("R6").Value = 13695.66855334834
$ws.Range("S6").Value = 0.06705461472320938
$ws.Range("T6").Value = 0.07512140717084029

# Row 7
$ws.Range("I7").Value = 0.3897411505765819
$ws.Range("J7").Value = 0.4097225989911443
$ws.Range("M7").Value = 2.005664333333333
$ws.Range("N7").Value = 6.016992999999999
$ws.Range("O7").Value = 0.02976781902817159
$ws.Range("P7").Value = 0.03172257287647481
$ws.Range("Q7").Value = 263.2905888958193
$ws.Range("R7").Value = 2369.615300062374
$ws.Range("S7").Value = 0.01160174403819506
$ws.Range("T7").Value = 0.01299745500563524

# Row 8
$ws.Range("I8").Value = 0.3897411505765819
$ws.Range("J8").Value = 0.4097225989911443
$ws.Range("O8").Value = 0.4846964599741412
$ws.Range("P8").Value = 0.5165248673390457
$ws.Range("S8").Value = 0.188906155990718
$ws.Range("T8").Value = 0.2116319110897098

# Row 9
$ws.Range("I9").Value = 0.3897411505765819
$ws.Range("J9").Value = 0.4097225989911443
$ws.Range("M9").Value = 10.495667
$ws.Range("N9").Value = 31.487001
$ws.Range("O9").Value = 0.1557753760903259
$ws.Range("P9").Value = 0.1660046278737794
$ws.Range("Q9").Value = 1377.803004898502
$ws.Range("R9").Value = 12400.22704408652
$ws.Range("S9").Value = 0.06071207430894337
$ws.Range("T9").Value = 0.06801584757700264

# Row 10
$ws.Range("I10").Value = 0.3897411505765819
$ws.Range("J10").Value = 0.4097225989911443
$ws.Range("M10").Value = 12.4553565
$ws.Range("N10").Value = 24.910713
$ws.Range("O10").Value = 0.1848608423958749
$ws.Range("P10").Value = 0.1313333601264699
$ws.Range("Q10").Value = 1635.058316234889
$ws.Range("R10").Value = 9810.349897409335
$ws.Range("S10").Value = 0.07204787741192445
$ws.Range("T10").Value = 0.05381024564525717

# Row 11
$ws.Range("I11").Value = 0.3897411505765819
$ws.Range("J11").Value = 0.4097225989911443
$ws.Range("M11").Value = 9.762884
$ws.Range("N11").Value = 29.288652
$ws.Range("O11").Value = 0.1448995025114864
$ws.Range("P11").Value = 0.1544145717842301
$ws.Range("Q11").Value = 1281.608011351304
$ws.Range("R11").Value = 11534.47210216174
$ws.Range("S11").Value = 0.05647329882680102
$ws.Range("T11").Value = 0.06326713967353936

# Row 12
$ws.Range("G12").Value = 0.2461213333333333
$ws.Range("H12").Value = 0.738364
$ws.Range("I12").Value = 0.0007307156985262189
$ws.Range("J12").Value = 0.0007681784042585035
$ws.Range("M12").Value = 2.005664333333333
$ws.Range("N12").Value = 6.016992999999999
$ws.Range("O12").Value = 0.02976781902817159
$ws.Range("P12").Value = 0.03172257287647481
$ws.Range("Q12").Value = 0.4936367799391111
$ws.Range("R12").Value = 4.442731019451999
$ws.Range("S12").Value = 0.00002175181267477247
$ws.Range("T12").Value = 0.00002436859541122451

# Row 13
$ws.Range("G13").Value = 0.2461213333333333
$ws.Range("H13").Value = 0.738364
$ws.Range("I13").Value = 0.0007307156985262189
$ws.Range("J13").Value = 0.0007681784042585035
$ws.Range("O13").Value = 0.4846964599741412
$ws.Range("P13").Value = 0.5165248673390457
$ws.Range("Q13").Value = 8.037673150427556
$ws.Range("R13").Value = 72.339058353848
$ws.Range("S13").Value = 0.0003541753123231901
$ws.Range("T13").Value = 0.0003967832483523434

# Row 14
$ws.Range("G14").Value = 0.2461213333333333
$ws.Range("H14").Value = 0.738364
$ws.Range("I14").Value = 0.0007307156985262189
$ws.Range("J14").Value = 0.0007681784042585035
$ws.Range("M14").Value = 10.495667
$ws.Range("N14").Value = 31.487001
$ws.Range("O14").Value = 0.1557753760903259
$ws.Range("P14").Value = 0.1660046278737794
$ws.Range("Q14").Value = 2.583207556262667
$ws.Range("R14").Value = 23.248868006364
$ws.Range("S14").Value = 0.0001138275127530269
$ws.Range("T14").Value = 0.0001275211701396065

# Row 15
$ws.Range("G15").Value = 0.2461213333333333
$ws.Range("H15").Value = 0.738364
$ws.Range("I15").Value = 0.0007307156985262189
$ws.Range("J15").Value = 0.0007681784042585035
$ws.Range("M15").Value = 12.4553565
$ws.Range("N15").Value = 24.910713
$ws.Range("O15").Value = 0.1848608423958749
$ws.Range("P15").Value = 0.1313333601264699
$ws.Range("Q15").Value = 3.065528948922
$ws.Range("R15").Value = 18.393173693532
$ws.Range("S15").Value = 0.000135080719581447
$ws.Range("T15").Value = 0.000100887451007859

# Row 16
$ws.Range("G16").Value = 0.2461213333333333
$ws.Range("H16").Value = 0.738364
$ws.Range("I16").Value = 0.0007307156985262189
$ws.Range("J16").Value = 0.0007681784042585035
$ws.Range("M16").Value = 9.762884
$ws.Range("N16").Value = 29.288652
$ws.Range("O16").Value = 0.1448995025114864
$ws.Range("P16").Value = 0.1544145717842301
$ws.Range("Q16").Value = 2.402854027258666
$ws.Range("R16").Value = 21.625686245328
$ws.Range("S16").Value = 0.0001058803411937824
$ws.Range("T16").Value = 0.00011861793934747

# Row 17
$ws.Range("G17").Value = 49.2786865
$ws.Range("H17").Value = 98.557373
$ws.Range("I17").Value = 0.1463047080910041
$ws.Range("J17").Value = 0.1025370217386683
$ws.Range("M17").Value = 2.005664333333333
$ws.Range("N17").Value = 6.016992999999999
$ws.Range("O17").Value = 0.02976781902817159
$ws.Range("P17").Value = 0.03172257287647481
$ws.Range("Q17").Value = 98.83650390656483
$ws.Range("R17").Value = 593.0190234393889
$ws.Range("S17").Value = 0.004355172073422482
$ws.Range("T17").Value = 0.003252738144641589

# Row 18
$ws.Range("G18").Value = 49.2786865
$ws.Range("H18").Value = 98.557373
$ws.Range("I18").Value = 0.1463047080910041
$ws.Range("J18").Value = 0.1025370217386683
$ws.Range("O18").Value = 0.4846964599741412
$ws.Range("P18").Value = 0.5165248673390457
$ws.Range("Q18").Value = 1609.311838210098
$ws.Range("R18").Value = 9655.871029260587
$ws.Range("S18").Value = 0.07091337408925978
$ws.Range("T18").Value = 0.05296292155090651

# Row 19
$ws.Range("G19").Value = 49.2786865
$ws.Range("H19").Value = 98.557373
$ws.Range("I19").Value = 0.1463047080910041
$ws.Range("J19").Value = 0.1025370217386683
$ws.Range("M19").Value = 10.495667
$ws.Range("N19").Value = 31.487001
$ws.Range("O19").Value = 0.1557753760903259
$ws.Range("P19").Value = 0.1660046278737794
$ws.Range("Q19").Value = 517.2126837013956
$ws.Range("R19").Value = 3103.276102208373
$ws.Range("S19").Value = 0.0227906709266615
$ws.Range("T19").Value = 0.01702162013701327

# Row 20
$ws.Range("G20").Value = 49.2786865
$ws.Range("H20").Value = 98.557373
$ws.Range("I20").Value = 0.1463047080910041
$ws.Range("J20").Value = 0.1025370217386683
$ws.Range("M20").Value = 12.4553565
$ws.Range("N20").Value = 24.910713
$ws.Range("O20").Value = 0.1848608423958749
$ws.Range("P20").Value = 0.1313333601264699
$ws.Range("Q20").Value = 613.7836082092373
$ws.Range("R20").Value = 2455.134432836949
$ws.Range("S20").Value = 0.02704601158418559
$ws.Range("T20").Value = 0.0134665316023002

# Row 21
$ws.Range("G21").Value = 49.2786865
$ws.Range("H21").Value = 98.557373
$ws.Range("I21").Value = 0.1463047080910041
$ws.Range("J21").Value = 0.1025370217386683
$ws.Range("M21").Value = 9.762884
$ws.Range("N21").Value = 29.288652
$ws.Range("O21").Value = 0.1448995025114864
$ws.Range("P21").Value = 0.1544145717842301
$ws.Range("Q21").Value = 481.102099971866
$ws.Range("R21").Value = 2886.612599831196
$ws.Range("S21").Value = 0.02119947941747473
$ws.Range("T21").Value = 0.01583321030380676

# Row 22
$ws.Range("G22").Value = 0.1539513333333333
$ws.Range("H22").Value = 0.461854
$ws.Range("I22").Value = 0.0004570699116250634
$ws.Range("J22").Value = 0.0004805032053572586
$ws.Range("M22").Value = 2.005664333333333
$ws.Range("N22").Value = 6.016992999999999
$ws.Range("O22").Value = 0.02976781902817159
$ws.Range("P22").Value = 0.03172257287647481
$ws.Range("Q22").Value = 0.3087746983357778
$ws.Range("R22").Value = 2.778972285022
$ws.Range("S22").Value = 0.00001360597441247727
$ws.Range("T22").Value = 0.00001524279794932538

# Row 23
$ws.Range("G23").Value = 0.1539513333333333
$ws.Range("H23").Value = 0.461854
$ws.Range("I23").Value = 0.0004570699116250634
$ws.Range("J23").Value = 0.0004805032053572586
$ws.Range("O23").Value = 0.4846964599741412
$ws.Range("P23").Value = 0.5165248673390457
$ws.Range("Q23").Value = 5.027644217780889
$ws.Range("R23").Value = 45.248797960028
$ws.Range("S23").Value = 0.0002215401681253618
$ws.Range("T23").Value = 0.0002481918544031442

# Row 24
$ws.Range("G24").Value = 0.1539513333333333
$ws.Range("H24").Value = 0.461854
$ws.Range("I24").Value = 0.0004570699116250634
$ws.Range("J24").Value = 0.0004805032053572586
$ws.Range("M24").Value = 10.495667
$ws.Range("N24").Value = 31.487001
$ws.Range("O24").Value = 0.1557753760903259
$ws.Range("P24").Value = 0.1660046278737794
$ws.Range("Q24").Value = 1.615821928872667
$ws.Range("R24").Value = 14.542397359854
$ws.Range("S24").Value = 0.00007120023738296624
$ws.Range("T24").Value = 0.00007976575579748991

# Row 25
$ws.Range("G25").Value = 0.1539513333333333
$ws.Range("H25").Value = 0.461854
$ws.Range("I25").Value = 0.0004570699116250634
$ws.Range("J25").Value = 0.0004805032053572586
$ws.Range("M25").Value = 12.4553565
$ws.Range("N25").Value = 24.910713
$ws.Range("O25").Value = 0.1848608423958749
$ws.Range("P25").Value = 0.1313333601264699
$ws.Range("Q25").Value = 1.917518740317
$ws.Range("R25").Value = 11.505112441902
$ws.Range("S25").Value = 0.00008449432889681729
$ws.Range("T25").Value = 0.00006310610051110796

# Row 26
$ws.Range("G26").Value = 0.1539513333333333
$ws.Range("H26").Value = 0.461854
$ws.Range("I26").Value = 0.0004570699116250634
$ws.Range("J26").Value = 0.0004805032053572586
$ws.Range("M26").Value = 9.762884
$ws.Range("N26").Value = 29.288652
$ws.Range("O26").Value = 0.1448995025114864
$ws.Range("P26").Value = 0.1544145717842301
$ws.Range("Q26").Value = 1.503009008978667
$ws.Range("R26").Value = 13.527081080808
$ws.Range("S26").Value = 0.00006622920280744071
$ws.Range("T26").Value = 0.00007419669669619105
